$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.660.36"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.690.17"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.48"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3942"
$ws.Range("E7").Value = "  -0.62%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4057"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.489"
$ws.Range("E9").Value = "  -1.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.002"
$ws.Range("E10").Value = "  +0.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.51"
$ws.Range("E11").Value = "  -2.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08853"
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.254"
$ws.Range("E13").Value = "  -1.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.57"
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.048"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001315"
$ws.Range("E16").Value = "  -0.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.691.01"
$ws.Range("E17").Value = "  -0.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "99.61"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("E19").Value = "  -1.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.54"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.997"
$ws.Range("E21").Value = "  +3.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.006"
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.33"
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.649.72"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.306"
$ws.Range("E25").Value = "  +9.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.356"
$ws.Range("E26").Value = "  +1.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.76"
$ws.Range("E27").Value = "  +1.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.50"
$ws.Range("E28").Value = "  +2.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "135.74"
$ws.Range("E29").Value = "  +1.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.170"
$ws.Range("E30").Value = "  +0.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.635"
$ws.Range("E31").Value = "  +2.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.878.20"
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.08557"
$ws.Range("E33").Value = "  -1.37%  "
$ws.Range("E34").Value = "  -2.63%  "
$ws.Range("E35").Value = "  -3.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.24"
$ws.Range("E36").Value = "  +1.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2733"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.884"
$ws.Range("E38").Value = "  -4.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.49"
$ws.Range("E39").Value = "  -1.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09202"
$ws.Range("E40").Value = "  +2.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02724"
$ws.Range("E41").Value = "  -1.91%  "
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7640"
$ws.Range("E43").Value = "  -0.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.01"
$ws.Range("E44").Value = "  +2.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.589"
$ws.Range("E45").Value = "  +5.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7137"
$ws.Range("E46").Value = "  -0.64%  "
$ws.Range("E47").Value = "  +0.74%  "
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.318"
$ws.Range("E50").Value = "  +1.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07967"
$ws.Range("E51").Value = "  -0.33%  "
